# Deploy refresh: Alvearie -> LinuxForHealth rebrand + version/date bump
# for StructureDefinition-insured-category.xlsx

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/insured-category"
# Version
$meta.Range("B3").Value = "8.0.0"
# Date
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Binding Value Set URL (row 7 = Extension.value[x])
$elements.Range("Y7").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/insured-category"

# Column widened to fit the longer value-set URL text
$elements.Columns.Item(25).ColumnWidth = 55.5

# The root "Extension" row no longer carries the ele-1/ext-1 constraint
# text in the Constraint(s) column (it now only appears on Extension.extension)
$elements.Range("AI2").Value = ""
